# Rotate the data rows 4, 5, 6:
#   new row4 = old row5, new row5 = old row6, new row6 = old row4
#
# Every data cell in this sheet is stored as text (inline string), even
# values that look numeric or date-like (e.g. "0", "2022", "2022-05-01").
# Writing such strings back through COM makes Excel auto-detect them as
# numbers/dates unless the cell is explicitly formatted as Text first.
# So for every cell we: force a Text number format, assign the string,
# then restore the cell's original (unstyled/"General") look by copying
# the style from a same-column cell that was never touched (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 17  # columns A..Q

# Capture original values for rows 4, 5, 6 (columns A..Q) before overwriting
$row4 = @()
$row5 = @()
$row6 = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $row4 += $ws.Cells.Item(4, $c).Value2
    $row5 += $ws.Cells.Item(5, $c).Value2
    $row6 += $ws.Cells.Item(6, $c).Value2
}

# Reference "plain"/unstyled style per column, taken from row 2 (untouched).
$plainStyle = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $plainStyle += $ws.Cells.Item(2, $c).Style
}

function Set-RowValues($targetRow, $values) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($targetRow, $c)
        $cell.NumberFormat = "@"
        $cell.Value2 = $values[$c - 1]
        $cell.Style = $plainStyle[$c - 1]
    }
}

# Write new row 4 = old row 5
Set-RowValues 4 $row5

# Write new row 5 = old row 6
Set-RowValues 5 $row6

# Write new row 6 = old row 4
Set-RowValues 6 $row4
